$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking.com crypto-price refresh (GitHub Actions scheduled update).
# Column D holds price strings exactly as scraped (dot-grouped thousands
# like "42.875.84", or plain decimals like "8.40"). Prefixing the plain
# decimals with an apostrophe forces Excel to store them as literal text
# instead of auto-parsing them into numbers (which would silently drop
# a significant trailing zero, e.g. "8.40" -> 8.4).
$quote = "'"

$ws.Range("D2").Value = '42.875.84'
$ws.Range("E2").Value = '  -0.32%  '

$ws.Range("D3").Value = '2.385.41'
$ws.Range("E3").Value = '  +3.56%  '

$ws.Range("E4").Value = '  -0.73%  '

$ws.Range("D5").Value = $quote + '328.32'
$ws.Range("E5").Value = '  +6.36%  '

$ws.Range("D6").Value = $quote + '99.57'
$ws.Range("E6").Value = '  -7.55%  '

$ws.Range("D7").Value = $quote + '0.636'
$ws.Range("E7").Value = '  +0.56%  '

$ws.Range("E8").Value = '  -0.21%  '

$ws.Range("D9").Value = $quote + '0.626'
$ws.Range("E9").Value = '  +2.14%  '

$ws.Range("D10").Value = $quote + '39.91'
$ws.Range("E10").Value = '  -8.46%  '

$ws.Range("D11").Value = $quote + '0.0922'
$ws.Range("E11").Value = '  -0.32%  '

$ws.Range("D12").Value = $quote + '8.40'
$ws.Range("E12").Value = '  -3.80%  '

$ws.Range("E13").Value = '  -2.43%  '

$ws.Range("D14").Value = $quote + '0.105'
$ws.Range("E14").Value = '  +0.94%  '

$ws.Range("E15").Value = '  +5.18%  '

$ws.Range("D16").Value = '2.734.44'
$ws.Range("E16").Value = '  +2.95%  '

$ws.Range("D17").Value = '2.384.37'
$ws.Range("E17").Value = '  -0.86%  '

$ws.Range("D18").Value = '42.782.34'
$ws.Range("E18").Value = '  -0.60%  '

$ws.Range("D19").Value = $quote + '7.77'
$ws.Range("E19").Value = '  +7.95%  '

$ws.Range("E20").Value = '  -1.16%  '

$ws.Range("D21").Value = $quote + '3.77'
$ws.Range("E21").Value = '  +8.49%  '

$ws.Range("D22").Value = $quote + '75.35'
$ws.Range("E22").Value = '  -0.42%  '

$ws.Range("D23").Value = $quote + '271.93'
$ws.Range("E23").Value = '  +6.74%  '

$ws.Range("D24").Value = $quote + '2.34'
$ws.Range("E24").Value = '  -6.55%  '

$ws.Range("E25").Value = '  +13.58%  '

$ws.Range("D26").Value = $quote + '0.998'
$ws.Range("E26").Value = '  -0.24%  '

$ws.Range("D27").Value = $quote + '11.51'
$ws.Range("E27").Value = '  -2.23%  '

$ws.Range("D28").Value = $quote + '23.91'
$ws.Range("E28").Value = '  +7.24%  '

$ws.Range("E29").Value = '  -1.69%  '

$ws.Range("D30").Value = $quote + '173.05'
$ws.Range("E30").Value = '  -0.09%  '

$ws.Range("E31").Value = '  -1.45%  '

$ws.Range("D32").Value = $quote + '0.0906'
$ws.Range("E32").Value = '  +0.97%  '

$ws.Range("D33").Value = $quote + '35.35'
$ws.Range("E33").Value = '  -7.95%  '

$ws.Range("E34").Value = '  +3.61%  '

$ws.Range("E35").Value = '  +2.24%  '

$ws.Range("D36").Value = $quote + '4.63'
$ws.Range("E36").Value = '  -7.74%  '

$ws.Range("E37").Value = '  -4.42%  '

$ws.Range("D38").Value = $quote + '3.86'
$ws.Range("E38").Value = '  -6.01%  '

$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = $quote + '2.88'
$ws.Range("E39").Value = '  +9.65%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = $quote + '0.105'
$ws.Range("E40").Value = '  +1.40%  '

$ws.Range("E41").Value = '  +2.95%  '

$ws.Range("E42").Value = '  -1.60%  '

$ws.Range("D43").Value = $quote + '68.93'
$ws.Range("E43").Value = '  -3.20%  '

$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").Value = $quote + '94.39'
$ws.Range("E44").Value = '  +49.27%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = $quote + '1.00'
$ws.Range("E45").Value = '  -0.09%  '

$ws.Range("D46").Value = $quote + '116.45'
$ws.Range("E46").Value = '  +7.48%  '

$ws.Range("D47").Value = $quote + '11.85'
$ws.Range("E47").Value = '  -3.60%  '

$ws.Range("D48").Value = $quote + '5.44'
$ws.Range("E48").Value = '  -5.08%  '

$ws.Range("E49").Value = '  +0.15%  '

$ws.Range("D50").Value = '1.623.15'
$ws.Range("E50").Value = '  +9.89%  '

$ws.Range("E51").Value = '  -1.75%  '
